$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '22.477.76'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.32%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.576.38'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1.000'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.09%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '288.72'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.67%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3691'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.92%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '47.80'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.79%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3331'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.71%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.149'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.85%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07562'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.02%  '
$ws.Range('E12').Value = '  -0.03%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.80'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.51%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.960'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.36%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.960'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.14%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.569.07'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.24%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001124'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.06%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '87.94'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.28%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06740'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.08%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.401'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.63%  '
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.000'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.12%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '16.56'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.04%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.03'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.49%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '22.477.64'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.37%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.389'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.76%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.652'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.56%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '150.32'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.44%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.68'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.26%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.995'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.34%  '
$ws.Range('E30').Value = '  +1.60%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.748.81'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.47%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.088'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.77%  '
$ws.Range('E33').Value = '  +0.11%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.993'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.26%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.900'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.06%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.08375'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.37%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02471'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2241'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.96%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06409'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.06%  '
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.297'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.01%  '
$ws.Range('B41').Value = 'InternetComputer(DFINITY)'
$ws.Range('C41').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.369'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.59%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.51'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.86%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.6291'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.54%  '
$ws.Range('B44').Value = 'Frax'
$ws.Range('C44').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.000'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.06%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '14.05'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.12%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6129'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +6.75%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.782'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.34%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.062'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.24%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '125.70'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.69%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.213'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.44%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07229'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.12%  '
